$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) sometimes holds numeric-looking text (e.g. "1.00",
# "5.30") that must stay text -- force text format before assigning so Excel
# does not silently coerce it to a Number and strip formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.746.87"
$ws.Range("E2").Value = "  +7.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.856.08"
$ws.Range("E3").Value = "  +10.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "424.43"
$ws.Range("E5").Value = "  +7.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.45"
$ws.Range("E6").Value = "  +3.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.842.67"
$ws.Range("E7").Value = "  +10.59%  "

$ws.Range("E8").Value = "  +2.48%  "

$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("E10").Value = "  +5.43%  "

$ws.Range("E11").Value = "  +5.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000344"
$ws.Range("E12").Value = "  +16.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.97"
$ws.Range("E13").Value = "  +2.24%  "

$ws.Range("E14").Value = "  +9.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.474.74"
$ws.Range("E15").Value = "  +10.20%  "

$ws.Range("E16").Value = "  +25.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.869.35"
$ws.Range("E17").Value = "  +10.00%  "

$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.01"
$ws.Range("E19").Value = "  +4.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.024.84"
$ws.Range("E20").Value = "  +7.40%  "

$ws.Range("E21").Value = "  +4.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "413.39"
$ws.Range("E22").Value = "  +2.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.98"
$ws.Range("E23").Value = "  +4.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.42"
$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("E25").Value = "  +5.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.71"
$ws.Range("E26").Value = "  +11.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("E27").Value = "  +10.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.30"
$ws.Range("E29").Value = "  +3.61%  "

$ws.Range("E30").Value = "  +37.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "725.00"
$ws.Range("E31").Value = "  +10.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.19"
$ws.Range("E32").Value = "  +9.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  +9.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.78"
$ws.Range("E34").Value = "  +5.29%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.17"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.48"
$ws.Range("E38").Value = "  +1.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.44"
$ws.Range("E39").Value = "  +31.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0753"
$ws.Range("E40").Value = "  +26.32%  "

$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  +4.91%  "

$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.135"
$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("E45").Value = "  +6.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.12"
$ws.Range("E46").Value = "  +4.05%  "

$ws.Range("E47").Value = "  +13.30%  "

$ws.Range("E48").Value = "  +4.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.17"
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("E50").Value = "  +3.95%  "

$ws.Range("E51").Value = "  +2.74%  "
